$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for the first file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 07:11:26"

# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime for the first file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 07:11:22"
$wsZhCn.Range("K2").Value = "2016-08-28 07:11:38"

# Sheet "de-de": Correspond Handoff Datetime / Correspond Handback DateTime for the first file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 07:11:26"
$wsDeDe.Range("K2").Value = "2016-08-28 07:11:45"
